$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "Hôtel des Andelys"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "US$1,423"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "4.0"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "Review score"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "929"

$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "City Inn Paris"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "US$608"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "6.7"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "Review score"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2,498"

$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "Austin's Saint Lazare Hotel"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "US$3,534"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "8.1"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "Very Good"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2,540"

$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "Auriane Porte De Versailles"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "US$1,948"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "5.6"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "Review score"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1,636"

$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "Zoku Paris"
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "US$3,315"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "9.1"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "Wonderful"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "887"

$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "Hotel Royal Phare"
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "US$3,701"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "8.3"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "Very Good"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1,735"

$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "Hôtel La Conversation"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "US$4,005"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "9.0"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "Wonderful"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "91"

$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = "Hôtel de l'Aveyron"
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "US$2,075"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "6.1"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "Review score"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "2,077"

$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "Hotel Anya"
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "US$1,918"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "7.4"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "Good"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1,132"

$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = "Glasgow Monceau by Patrick Hayat"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "US$2,897"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "8.0"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "Very Good"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "674"

$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = "Printania Porte de Versailles"
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "US$1,948"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "5.7"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "Review score"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1,284"

$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = "PORTE MAILLOT CHAMPS ELYSÉES PARIS"
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "US$3,469"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "4.7"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "Review score"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "3"

$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = "STYLE HOTEL"
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "US$1,939"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "6.0"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "Review score"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "744"

$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = "Austin's Arts Et Metiers Hotel"
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "US$3,971"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "8.2"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "Very Good"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "2,030"

$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = "Hôtel De Castiglione"
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "US$4,559"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "7.5"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "Good"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "3,693"

$ws.Range("A17").NumberFormat = "@"
$ws.Range("A17").Value = "FM Hotel"
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "US$1,570"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "3.6"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "Review score"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "853"

$ws.Range("A18").NumberFormat = "@"
$ws.Range("A18").Value = "Charmant studio grand balcon parking"
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "US$2,179"
$ws.Range("C18").Value = ""
$ws.Range("D18").Value = ""
$ws.Range("E18").Value = ""

$ws.Range("A19").NumberFormat = "@"
$ws.Range("A19").Value = "Aparthotel Adagio Paris Buttes Chaumont"
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "US$4,354"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "8.3"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "Very Good"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "562"

$ws.Range("A20").NumberFormat = "@"
$ws.Range("A20").Value = "Apart hotel Le Marais Centre de Paris"
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "US$4,613"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "8.0"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "Very Good"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "163"

$ws.Range("A21").NumberFormat = "@"
$ws.Range("A21").Value = "Studio - Butte aux cailles 13eme"
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "US$2,034"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "7.7"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "Good"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "9"

$ws.Range("A22").NumberFormat = "@"
$ws.Range("A22").Value = "HOTEL SAVOY"
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "US$1,575"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "6.1"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "Review score"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "779"

$ws.Range("A23").NumberFormat = "@"
$ws.Range("A23").Value = "Hôtel Des Fontaines"
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "US$2,539"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "5.9"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "Review score"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "473"

$ws.Range("A24").NumberFormat = "@"
$ws.Range("A24").Value = "Hotel Eden Montmartre"
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "US$2,438"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "7.9"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "Good"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1,595"

$ws.Range("A25").NumberFormat = "@"
$ws.Range("A25").Value = "Domitys L'Ellipse"
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = "US$2,977"
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = "9.2"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "Wonderful"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "6"

$ws.Range("A26").NumberFormat = "@"
$ws.Range("A26").Value = "Hotel de France 18"
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "US$1,279"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "6.0"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "Review score"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "1,891"
